$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (J1, K1) ---
$ws.Range("J1").Value = "id"
$ws.Range("K1").Value = "link"

# --- Remove the old guest rows (2:4); data is being relocated to rows 6:9 ---
$ws.Rows("2:4").ClearContents()

# --- Row 6: גיא מדואל ---
$ws.Range("A6").Value = "גיא מדואל"
$ws.Range("B6").Value = "guymaduel2302@gmail.com"
$ws.Range("C6").Value = "0522505756"
$ws.Range("D6").Value = "Israel"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("I6").Value = $false

# --- Row 7: Aaliyah Hay ---
$ws.Range("A7").Value = "Aaliyah Hay"
$ws.Range("B7").Value = "aaliyah.a.hay@gmail.com"
$ws.Range("C7").Value = "0505505756"
$ws.Range("D7").Value = "America"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = $true

# --- Row 8: guy maduel ---
$ws.Range("A8").Value = "guy maduel"
$ws.Range("B8").Value = "guy.gm.maduel@gmail.com"
$ws.Range("C8").Value = "0522505756"
$ws.Range("D8").Value = "America"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("I8").Value = $true

# --- Row 9: Aaliyah Hay / whatsapp conflict row ---
$ws.Range("A9").Value = "Aaliyah Hay"
$ws.Range("B9").Value = "whatsapp"
$ws.Range("C9").Value = "0505505756"
$ws.Range("D9").Value = "America"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0
$ws.Range("I9").Value = $false

# --- Column widths for B, C, J, K (closest achievable values given the
#     host's pixel-quantized ColumnWidth conversion) ---
$ws.Columns("B").ColumnWidth = 18
$ws.Columns("C").ColumnWidth = 16.3333333333333
$ws.Columns("J").ColumnWidth = 1.83333333333333
$ws.Columns("K").ColumnWidth = 80.3333333333333

# --- Active selection ---
$ws.Range("B3").Select()
